# Extend the linear piecewise heat rates output-parameter workbook with the
# new bus-detailed / PTDF related cost outputs.
#
# Sheet layout (tab order):
#   1 = obj_output          (A: obj_class name, B: obj_name)
#   2 = obj_report          (unchanged data)
#   3 = rel_report__output  (A: rel_class_name, B: report, C: output)
#
# Strategy: clear the data area of sheet1 and sheet3, then rewrite every
# cell (values + styles) from scratch in the exact order required so that
# the shared-strings table grows with the new strings in the same order
# they were introduced upstream, and cell formatting (the "quote prefix"
# styles used for a couple of rows) is reproduced faithfully by copying
# format from a template cell that already carries that style.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # obj_output
$ws2 = $wb.Worksheets.Item(2)   # obj_report (left untouched)
$ws3 = $wb.Worksheets.Item(3)   # rel_report__output

# ---------------------------------------------------------------------
# Grab template cells (with their existing styles) from sheet1 BEFORE we
# clear anything, so we can re-apply the same two "quote-prefixed" cell
# styles (s="1" and s="2") to the rewritten cells further down. Parked on
# sheet2 (obj_report), far to the right of its used range, since sheet1 and
# sheet3 are about to be wiped and rebuilt.
# ---------------------------------------------------------------------
$styleA = $ws1.Range("B2")   # style s="2"
$styleB = $ws1.Range("B4")   # style s="1"

$styleA.Copy()
$ws2.Range("Z1").PasteSpecial(-4122)
$styleB.Copy()
$ws2.Range("Z2").PasteSpecial(-4122)
$tmplS2 = $ws2.Range("Z1")
$tmplS1 = $ws2.Range("Z2")

# Now clear both sheets completely and start fresh.
$ws1.Cells.Clear()
$ws3.Cells.Clear()

# ---------------------------------------------------------------------
# obj_output (sheet1) header
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "obj_class name"
$ws1.Range("B1").Value = "obj_name"

# Data rows 2-25. Column A is always "output".
$ws1.Range("A2:A25").Value = "output"

$ws1.Range("B2").Value = "unit_flow"
$ws1.Range("B3").Value = "units_started_up"
$ws1.Range("B4").Value = "units_shut_down"
$ws1.Range("B5").Value = "units_available"
$ws1.Range("B6").Value = "units_on"
$ws1.Range("B7").Value = "connection_flow"
$ws1.Range("B8").Value = "node_state"
$ws1.Range("B9").Value = "node_slack_pos"
$ws1.Range("B10").Value = "node_slack_neg"
$ws1.Range("B11").Value = "node_injection"
$ws1.Range("B12").Value = "total_costs"
$ws1.Range("B13").Value = "start_up_costs"
$ws1.Range("B14").Value = "shut_down_costs"
$ws1.Range("B15").Value = "fuel_costs"
$ws1.Range("B16").Value = "operating_costs"
$ws1.Range("B17").Value = "fixed_om_costs"
$ws1.Range("B18").Value = "variable_om_costs"
$ws1.Range("B19").Value = "ramp_costs"
$ws1.Range("B20").Value = "res_proc_costs"
$ws1.Range("B21").Value = "renewable_curtailment_costs"
# NOTE: "taxes" must be registered as a shared string before
# "connection_flow_costs" to match the upstream ordering, even though the
# "connection_flow_costs" row sits above the "taxes" row on the sheet.
$ws1.Range("B23").Value = "taxes"
$ws1.Range("B22").Value = "connection_flow_costs"
$ws1.Range("B24").Value = "investment_costs"
$ws1.Range("B25").Value = "objective_penalties"

# Re-apply the two special cell styles.
$tmplS2.Copy()
$ws1.Range("B2").PasteSpecial(-4122)
$tmplS1.Copy()
$ws1.Range("B4").PasteSpecial(-4122)
$tmplS1.Copy()
$ws1.Range("B9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# rel_report__output (sheet3)
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "rel_class_name"
$ws3.Range("B1").Value = "report"
$ws3.Range("C1").Value = "output"

$ws3.Range("A2:A25").Value = "report__output"
$ws3.Range("B2:B25").Value = "result_temp"

$ws3.Range("C2").Value = "unit_flow"
$ws3.Range("C3").Value = "units_started_up"
$ws3.Range("C4").Value = "units_shut_down"
$ws3.Range("C5").Value = "units_available"
$ws3.Range("C6").Value = "units_on"
$ws3.Range("C7").Value = "connection_flow"
$ws3.Range("C8").Value = "node_state"
$ws3.Range("C9").Value = "node_slack_pos"
$ws3.Range("C10").Value = "node_slack_neg"
$ws3.Range("C11").Value = "node_injection"
$ws3.Range("C12").Value = "total_costs"
$ws3.Range("C13").Value = "start_up_costs"
$ws3.Range("C14").Value = "shut_down_costs"
$ws3.Range("C15").Value = "fuel_costs"
$ws3.Range("C16").Value = "operating_costs"
$ws3.Range("C17").Value = "fixed_om_costs"
$ws3.Range("C18").Value = "variable_om_costs"
$ws3.Range("C19").Value = "ramp_costs"
$ws3.Range("C20").Value = "res_proc_costs"
$ws3.Range("C21").Value = "renewable_curtailment_costs"
$ws3.Range("C22").Value = "connection_flow_costs"
$ws3.Range("C23").Value = "taxes"
$ws3.Range("C24").Value = "investment_costs"
$ws3.Range("C25").Value = "objective_penalties"

# Style: column B is entirely the s="1" style; column C keeps General
# style except for the same two rows that are special-styled on sheet1
# (shifted down by one row since sheet3 has the extra leading column).
$tmplS1.Copy()
$ws3.Range("B2:B25").PasteSpecial(-4122)

$tmplS2.Copy()
$ws3.Range("C2").PasteSpecial(-4122)
$tmplS1.Copy()
$ws3.Range("C4").PasteSpecial(-4122)
$tmplS1.Copy()
$ws3.Range("C9").PasteSpecial(-4122)

# Trailing placeholder rows (26-35), column B only, alternating style.
$tmplS1.Copy()
$ws3.Range("B26:B27").PasteSpecial(-4122)
$tmplS2.Copy()
$ws3.Range("B28:B35").PasteSpecial(-4122)

# Drop the scratch template cells used to carry styles around.
$ws2.Range("Z1:Z2").Clear()

# ---------------------------------------------------------------------
# Selections / active sheet, matching the upstream commit: the active
# workbook tab moves from obj_output to rel_report__output, and both
# sheets get new selections reflecting the review of the new cost rows.
# ---------------------------------------------------------------------
$ws1.Range("B2:B25").Select()
$ws3.Range("E19").Select()
$ws3.Activate()
